$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column: header in H1 (reusing the same header style as the
# existing G1 "sum" header via copy/paste-formats so it lands on the same
# style index), and the data value 1 in H2.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H2").Value = 1
